# Daily attendance processing - 2025-12-07 15:25:07
# Normalize the "Recorded By" column (G): when the value is a comma-separated
# list that starts with the token "System", move that leading "System" token
# to the end of the list instead of the front.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -and $value.StartsWith("System, ")) {
        $parts = $value.Split(",")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }
        $first = $parts[0]
        $rest = $parts[1..($parts.Length - 1)]
        $newParts = $rest + @($first)
        $newValue = [string]::Join(", ", $newParts)
        $cell.Value = $newValue
    }
}
